# Update "想去人数" (want-to-go count, column F) figures across all sheets
# to reflect newly generated output (gh-pages rebuild at 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 993
$ws1.Range("F7").Value = 1224
$ws1.Range("F9").Value = 45
$ws1.Range("F11").Value = 1061
$ws1.Range("F12").Value = 4552
$ws1.Range("F15").Value = 1736
$ws1.Range("F17").Value = 669
$ws1.Range("F18").Value = 27
$ws1.Range("F20").Value = 390
$ws1.Range("F21").Value = 1104
$ws1.Range("F22").Value = 1547
$ws1.Range("F24").Value = 692
$ws1.Range("F25").Value = 527
$ws1.Range("F27").Value = 622
$ws1.Range("F28").Value = 93
$ws1.Range("F30").Value = 1174
$ws1.Range("F31").Value = 358
$ws1.Range("F32").Value = 2482
$ws1.Range("F34").Value = 1474
$ws1.Range("F35").Value = 475
$ws1.Range("F38").Value = 4122

# Sheet 2: 演出
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F6").Value = 207
$ws2.Range("F23").Value = 269
$ws2.Range("F24").Value = 5
$ws2.Range("F25").Value = 5
$ws2.Range("F40").Value = 24

# Sheet 3: 本地生活
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F4").Value = 1307
$ws3.Range("F5").Value = 1706
$ws3.Range("F8").Value = 162

# Sheet 4: 全部类型
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1307
$ws4.Range("F3").Value = 1706
$ws4.Range("F7").Value = 993
$ws4.Range("F9").Value = 1224
$ws4.Range("F12").Value = 45
$ws4.Range("F14").Value = 207
$ws4.Range("F15").Value = 207
$ws4.Range("F16").Value = 162
$ws4.Range("F18").Value = 1061
$ws4.Range("F20").Value = 4552
$ws4.Range("F23").Value = 1736
$ws4.Range("F25").Value = 669
$ws4.Range("F27").Value = 390
$ws4.Range("F28").Value = 1104
$ws4.Range("F29").Value = 1547
$ws4.Range("F32").Value = 692
$ws4.Range("F33").Value = 527
$ws4.Range("F35").Value = 622
$ws4.Range("F36").Value = 93
$ws4.Range("F39").Value = 269
$ws4.Range("F40").Value = 5
$ws4.Range("F42").Value = 1174
$ws4.Range("F43").Value = 358
$ws4.Range("F44").Value = 2482
$ws4.Range("F46").Value = 1474
$ws4.Range("F47").Value = 475
$ws4.Range("F49").Value = 4122
$ws4.Range("F51").Value = 24
